$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.492.32'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '1.556.22'
$ws.Range("E3").Value = '  -1.51%  '
$ws.Range("E4").Value = '  -0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.87'
$ws.Range("E5").Value = '  -1.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.483'
$ws.Range("E6").Value = '  -1.70%  '
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.30'
$ws.Range("E8").Value = '  +1.89%  '
$ws.Range("E9").Value = '  -1.69%  '
$ws.Range("E10").Value = '  -0.91%  '
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("D12").Value = '1.778.13'
$ws.Range("E12").Value = '  -1.57%  '
$ws.Range("D13").Value = '1.569.41'
$ws.Range("E13").Value = '  -0.79%  '
$ws.Range("D14").Value = '28.485.25'
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("E15").Value = '  -1.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.511'
$ws.Range("E16").Value = '  -1.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.29'
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.63'
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.38'
$ws.Range("E19").Value = '  -0.83%  '
$ws.Range("E20").Value = '  -2.08%  '
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("E22").Value = '  -0.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.94'
$ws.Range("E23").Value = '  -1.71%  '
$ws.Range("E24").Value = '  +1.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.60'
$ws.Range("E25").Value = '  -0.74%  '
$ws.Range("E26").Value = '  -1.71%  '
$ws.Range("E27").Value = '  -0.86%  '
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("E29").Value = '  -2.54%  '
$ws.Range("E30").Value = '  -3.28%  '
$ws.Range("E31").Value = '  -3.99%  '
$ws.Range("E32").Value = '  -1.04%  '
$ws.Range("D33").Value = '1.396.31'
$ws.Range("E33").Value = '  -0.17%  '
$ws.Range("E34").Value = '  -2.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.05'
$ws.Range("E35").Value = '  -4.12%  '
$ws.Range("E36").Value = '  -1.34%  '
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.66'
$ws.Range("E37").Value = '  +0.38%  '
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.29'
$ws.Range("E38").Value = '  -3.07%  '
$ws.Range("E39").Value = '  -1.32%  '
$ws.Range("E40").Value = '  +2.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.519'
$ws.Range("E41").Value = '  -0.83%  '
$ws.Range("E42").Value = '  -0.39%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.777'
$ws.Range("E43").Value = '  -1.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0467'
$ws.Range("E44").Value = '  +2.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.38'
$ws.Range("E45").Value = '  +2.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.34'
$ws.Range("E46").Value = '  -1.53%  '
$ws.Range("D47").Value = '1.691.19'
$ws.Range("E47").Value = '  -1.58%  '
$ws.Range("E48").Value = '  -6.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '43.90'
$ws.Range("E49").Value = '  +7.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '85.49'
$ws.Range("E50").Value = '  -0.71%  '
$ws.Range("E51").Value = '  -1.57%  '
